# Update cryptocurrency price/volume data (scraped refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.302.44'
$ws.Range('E2').Value = '  -1.37%  '

$ws.Range('D3').Value = '2.359.23'
$ws.Range('E3').Value = '  +5.28%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '''232.44'
$ws.Range('E5').Value = '  +0.86%  '

$ws.Range('E6').Value = '  +0.88%  '

$ws.Range('D7').Value = '''68.40'
$ws.Range('E7').Value = '  +8.28%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').Value = '''0.458'
$ws.Range('E9').Value = '  +1.21%  '

$ws.Range('D10').Value = '''0.0952'
$ws.Range('E10').Value = '  -2.13%  '

$ws.Range('D11').Value = '''56.76'
$ws.Range('E11').Value = '  +0.10%  '

$ws.Range('D12').Value = '''26.48'
$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('D13').Value = '2.711.75'
$ws.Range('E13').Value = '  +5.05%  '

$ws.Range('E14').Value = '  -0.48%  '

$ws.Range('D15').Value = '''15.60'
$ws.Range('E15').Value = '  +1.02%  '

$ws.Range('E16').Value = '  +3.02%  '

$ws.Range('E17').Value = '  +1.86%  '

$ws.Range('D18').Value = '2.361.61'
$ws.Range('E18').Value = '  +5.28%  '

$ws.Range('D19').Value = '43.283.91'
$ws.Range('E19').Value = '  -0.94%  '

$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('E20').Value = '  -0.60%  '

$ws.Range('D21').Value = '''74.01'
$ws.Range('E21').Value = '  +1.39%  '

$ws.Range('D22').Value = '''6.26'
$ws.Range('E22').Value = '  +4.21%  '

$ws.Range('D23').Value = '''248.52'
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('D24').Value = '''4.04'
$ws.Range('E24').Value = '  +18.75%  '

$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('E26').Value = '  +1.14%  '

$ws.Range('E27').Value = '  +0.72%  '

$ws.Range('E28').Value = '  -3.90%  '

$ws.Range('E29').Value = '  +7.50%  '

$ws.Range('D30').Value = '''173.98'
$ws.Range('E30').Value = '  +2.08%  '

$ws.Range('E31').Value = '  +11.98%  '

$ws.Range('E32').Value = '  -6.53%  '

$ws.Range('E33').Value = '  +0.45%  '

$ws.Range('D34').Value = '''5.02'
$ws.Range('E34').Value = '  +5.56%  '

$ws.Range('E35').Value = '  -0.88%  '

$ws.Range('D36').Value = '''5.06'
$ws.Range('E36').Value = '  +4.27%  '

$ws.Range('D37').Value = '''2.52'
$ws.Range('E37').Value = '  +10.92%  '

$ws.Range('E38').Value = '  +1.58%  '

$ws.Range('D39').Value = '''3.65'
$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('D40').Value = '''0.0254'
$ws.Range('E40').Value = '  -1.50%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''8.95'
$ws.Range('E41').Value = '  +8.99%  '

$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.05%  '

$ws.Range('D43').Value = '''18.09'
$ws.Range('E43').Value = '  +5.20%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '''1.16'
$ws.Range('E44').Value = '  +8.27%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''98.76'
$ws.Range('E45').Value = '  +1.92%  '

$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '''1.21'
$ws.Range('E46').Value = '  +2.30%  '

$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = '''4.46'
$ws.Range('E47').Value = '  +2.52%  '

$ws.Range('D48').Value = '''0.0949'
$ws.Range('E48').Value = '  +0.82%  '

$ws.Range('D49').Value = '1.447.58'
$ws.Range('E49').Value = '  +1.42%  '

$ws.Range('D50').Value = '2.584.04'
$ws.Range('E50').Value = '  +5.35%  '

$ws.Range('D51').Value = '''2.26'
$ws.Range('E51').Value = '  -2.67%  '
